# Updated cryptos list values (price + volume%) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.213.09"
$ws.Range("E2").Value = "  -2.20%  "

# Row 3
$ws.Range("D3").Value = "1.820.97"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.05%  "

# Row 5
$ws.Range("D5").Value = "'314.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

# Row 6
$ws.Range("D6").Value = "'1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
$ws.Range("D7").Value = "'0.4259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.28%  "

# Row 8
$ws.Range("E8").Value = "  -2.82%  "

# Row 9
$ws.Range("D9").Value = "'0.07225"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "

# Row 10
$ws.Range("D10").Value = "'0.8606"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "

# Row 11
$ws.Range("D11").Value = "'20.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.05%  "

# Row 12
$ws.Range("D12").Value = "1.832.44"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13
$ws.Range("D13").Value = "'6.666"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.58%  "

# Row 14
$ws.Range("D14").Value = "'0.07135"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("D15").Value = "'5.301"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.32%  "

# Row 16
$ws.Range("D16").Value = "'88.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.25%  "

# Row 17
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("D18").Value = "'0.000008853"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "

# Row 19
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "

# Row 20
$ws.Range("D20").Value = "'15.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.72%  "

# Row 21
$ws.Range("D21").Value = "27.254.35"
$ws.Range("E21").Value = "  -2.00%  "

# Row 22
$ws.Range("D22").Value = "'5.137"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "

# Row 23
$ws.Range("D23").Value = "'10.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "

# Row 24
$ws.Range("D24").Value = "2.046.98"
$ws.Range("E24").Value = "  -1.71%  "

# Row 25
$ws.Range("E25").Value = "  -1.18%  "

# Row 26
$ws.Range("D26").Value = "'153.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.04%  "

# Row 27
$ws.Range("D27").Value = "'18.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

# Row 28
$ws.Range("D28").Value = "'2.110"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.90%  "

# Row 29
$ws.Range("D29").Value = "'5.223"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

# Row 30
$ws.Range("D30").Value = "'116.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.50%  "

# Row 31
$ws.Range("D31").Value = "'0.08893"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "

# Row 32
$ws.Range("E32").Value = "  -2.83%  "

# Row 33
$ws.Range("D33").Value = "'0.7560"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "

# Row 34
$ws.Range("D34").Value = "'4.433"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "

# Row 35
$ws.Range("D35").Value = "'2.815"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.59%  "

# Row 36
$ws.Range("D36").Value = "'1.007"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "

# Row 37
$ws.Range("D37").Value = "'1.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "

# Row 38
$ws.Range("D38").Value = "'0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "

# Row 39
$ws.Range("D39").Value = "'0.05275"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.154"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.45%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.868"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.24%  "

# Row 42
$ws.Range("D42").Value = "'0.1686"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "

# Row 43
$ws.Range("D43").Value = "'0.5033"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44
$ws.Range("D44").Value = "'8.611"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "

# Row 45
$ws.Range("D45").Value = "'10.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "

# Row 46
$ws.Range("D46").Value = "'106.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.99%  "

# Row 47
$ws.Range("D47").Value = "'0.4728"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "

# Row 48
$ws.Range("D48").Value = "'1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.17%  "

# Row 49
$ws.Range("D49").Value = "'0.06372"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "

# Row 50
$ws.Range("D50").Value = "'1.655"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.96%  "

# Row 51
$ws.Range("D51").Value = "'1.802"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.52%  "
